$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.242.60"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").Value = "2.496.68"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.30%  "

# Row 7
$ws.Range("E7").Value = "  -1.83%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0801"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "

# Row 15
$ws.Range("D15").Value = "2.887.79"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16
$ws.Range("D16").Value = "2.476.35"
$ws.Range("E16").Value = "  -1.68%  "

# Row 17
$ws.Range("E17").Value = "  -1.14%  "

# Row 18
$ws.Range("D18").Value = "48.047.73"
$ws.Range("E18").Value = "  +0.22%  "

# Row 19
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.95%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0931"
$ws.Range("E22").Value = "  -1.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

# Row 28
$ws.Range("E28").Value = "  -7.95%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.61%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "

# Row 34
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.69%  "

# Row 36
$ws.Range("E36").Value = "  -1.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "

# Row 40
$ws.Range("E40").Value = "  -0.84%  "

# Row 41
$ws.Range("E41").Value = "  -0.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "

# Row 44
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
$ws.Range("D45").Value = "1.990.03"
$ws.Range("E45").Value = "  -1.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.27%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.90%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.38%  "

# Row 50
$ws.Range("E50").Value = "  -1.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
